$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, [string]$val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "94.437.22"
$ws.Range("E2").Value = "  -3.30%  "
Set-TextValue $ws.Range("D3") "3.430.43"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue $ws.Range("D5") "237.60"
$ws.Range("E5").Value = "  -5.88%  "
Set-TextValue $ws.Range("D6") "642.83"
$ws.Range("E6").Value = "  -2.21%  "
Set-TextValue $ws.Range("D7") "1.45"
$ws.Range("E7").Value = "  +0.91%  "
Set-TextValue $ws.Range("D8") "0.407"
$ws.Range("E8").Value = "  -3.31%  "
Set-TextValue $ws.Range("D10") "0.974"
$ws.Range("E10").Value = "  -5.82%  "
Set-TextValue $ws.Range("D11") "3.427.60"
$ws.Range("E11").Value = "  +1.85%  "
Set-TextValue $ws.Range("D12") "0.199"
$ws.Range("E12").Value = "  -4.79%  "
Set-TextValue $ws.Range("D13") "41.62"
$ws.Range("E13").Value = "  +0.69%  "
Set-TextValue $ws.Range("D14") "6.20"
$ws.Range("E14").Value = "  +2.15%  "
Set-TextValue $ws.Range("D15") "94.169.68"
$ws.Range("E15").Value = "  -3.41%  "
Set-TextValue $ws.Range("D16") "4.078.52"
$ws.Range("E16").Value = "  +2.12%  "
Set-TextValue $ws.Range("D17") "0.0000252"
$ws.Range("E17").Value = "  -0.89%  "
Set-TextValue $ws.Range("D18") "8.33"
$ws.Range("E18").Value = "  -5.91%  "
Set-TextValue $ws.Range("D19") "3.424.27"
$ws.Range("E19").Value = "  +1.69%  "
Set-TextValue $ws.Range("D20") "17.52"
$ws.Range("E20").Value = "  -2.87%  "
Set-TextValue $ws.Range("D21") "11.52"
$ws.Range("E21").Value = "  +6.03%  "
Set-TextValue $ws.Range("D22") "0.502"
$ws.Range("E22").Value = "  -4.99%  "
Set-TextValue $ws.Range("D23") "498.19"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("E24").Value = "  -4.89%  "
Set-TextValue $ws.Range("D25") "0.0000193"
$ws.Range("E25").Value = "  -2.89%  "
Set-TextValue $ws.Range("D26") "6.50"
$ws.Range("E26").Value = "  -8.30%  "
Set-TextValue $ws.Range("D27") "93.84"
$ws.Range("E27").Value = "  +0.41%  "
Set-TextValue $ws.Range("D28") "11.98"
$ws.Range("E28").Value = "  -2.36%  "
Set-TextValue $ws.Range("D29") "3.616.91"
$ws.Range("E29").Value = "  +1.94%  "
Set-TextValue $ws.Range("D30") "11.69"
$ws.Range("E30").Value = "  +3.26%  "
Set-TextValue $ws.Range("D31") "0.998"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").Value = "  +8.27%  "
Set-TextValue $ws.Range("D33") "0.139"
$ws.Range("E33").Value = "  -1.31%  "
Set-TextValue $ws.Range("D34") "0.999"
$ws.Range("E34").Value = "  +0.26%  "
Set-TextValue $ws.Range("D35") "0.179"
$ws.Range("E35").Value = "  -3.97%  "
Set-TextValue $ws.Range("D36") "29.73"
$ws.Range("E36").Value = "  +3.91%  "
Set-TextValue $ws.Range("D37") "0.554"
$ws.Range("E37").Value = "  -0.92%  "
Set-TextValue $ws.Range("D38") "544.69"
$ws.Range("E38").Value = "  +3.81%  "
Set-TextValue $ws.Range("D39") "7.66"
$ws.Range("E39").Value = "  -4.13%  "
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("E41").Value = "  +0.01%  "
Set-TextValue $ws.Range("D42") "0.151"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("E43").Value = "  +6.23%  "
Set-TextValue $ws.Range("D44") "24.09"
$ws.Range("E44").Value = "  -1.23%  "
Set-TextValue $ws.Range("D45") "1.72"
$ws.Range("E45").Value = "  -0.16%  "
Set-TextValue $ws.Range("D49") "0.0409"
$ws.Range("E49").Value = "  -4.49%  "
Set-TextValue $ws.Range("D50") "2.19"
$ws.Range("E50").Value = "  -4.26%  "
Set-TextValue $ws.Range("D51") "54.46"
$ws.Range("E51").Value = "  -2.67%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D46") "5.58"
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("B47").Value = "MantraDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue $ws.Range("D47") "3.61"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D48") "3.34"
$ws.Range("E48").Value = "  +5.48%  "
